$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.018.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.34%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.832.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9982'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6264'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07610'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2927'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07713'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.837.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.11%  '

$ws.Range("E13").Value = '  -0.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6659'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001023'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +17.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.046'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.036.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '226.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9984'
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.180'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9994'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.484'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1375'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.487'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.104'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.014'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.44%  '

$ws.Range("E31").Value = '  -1.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05215'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.843'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7361'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.139'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.58%  '

$ws.Range("E36").Value = '  +1.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.241.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.758'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01785'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.341'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8942'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9992'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.49'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.982.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.44%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5102'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.71%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4034'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.72%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.893'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.645'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05746'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.83%  '

